# Apply the data updates described in the diff:
#  1. Update the "as of" date in the confidential disclaimer text (A7) from
#     2021-03-22 to 2021-03-23.
#  2. Update the Weight/Percent Change values in rows 2-4 (columns D and E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to edit, then restore protection after.
$ws.Unprotect()

# 1. Update the disclaimer text in cell A7.
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# 2. Update the numeric Weight / Percent Change figures.
$ws.Range("D2").Value = 0.8417817162810487
$ws.Range("E2").Value = -0.01385620915032693

$ws.Range("D3").Value = 0.1582182837189513
$ws.Range("E3").Value = -0.01910237388724034

$ws.Range("E4").Value = -0.01468624833110832

# Restore sheet protection (matches original settings).
$ws.Protect()
